$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated strikeout (K) values regenerated from save_data (Strike# -> K), rows 2-22
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
